$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.300.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.386.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.194"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.590"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "680.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.929.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.358.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.37%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.393.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.120"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.904"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "557.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("E34").Value = "  +6.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.669.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.37%  "
$ws.Range("E38").Value = "  +3.99%  "
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0699"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.340"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("E48").Value = "  +5.99%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.10%  "
